$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 721.2727
$ws.Range("I33").Value = 294.3
$ws.Range("K33").Value = 294.3
$ws.Range("M33").Value = -65.30000000000001

$ws.Range("H40").Value = 15496.733
$ws.Range("I40").Value = 2612.75
$ws.Range("K40").Value = 2612.75
$ws.Range("M40").Value = -2437.75

$ws.Range("H43").Value = 16049
$ws.Range("I43").Value = 18066.334
$ws.Range("K43").Value = 18066.334
$ws.Range("M43").Value = -17997.334

$ws.Range("H113").Value = 4734.3
$ws.Range("I113").Value = 4616.6665
$ws.Range("J113").Value = 4784.7144
$ws.Range("K113").Value = 4616.6665
$ws.Range("L113").Value = 4784.7144
$ws.Range("M113").Value = -1362.6665
$ws.Range("N113").Value = -11292.7144

$ws.Range("H116").Value = 4842.7144
$ws.Range("I116").Value = 4816.6665
$ws.Range("J116").Value = 4999
$ws.Range("K116").Value = 4816.6665
$ws.Range("L116").Value = 4999
$ws.Range("M116").Value = -1374.6665
$ws.Range("N116").Value = -11883

$ws.Range("H131").Value = 1818
$ws.Range("I131").Value = 1224.5
$ws.Range("K131").Value = 3673.5
$ws.Range("M131").Value = 1366.5

$ws.Range("H138").Value = 3008.152
$ws.Range("I138").Value = 1290.7916
$ws.Range("J138").Value = 4881.636
$ws.Range("K138").Value = 3872.3748
$ws.Range("L138").Value = 14644.908
$ws.Range("M138").Value = 1267.6252
$ws.Range("N138").Value = -24924.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3605.7693
$ws.Range("I45").Value = 2616.8667
$ws.Range("K45").Value = 2616.8667
$ws.Range("M45").Value = -2239.8667

$ws.Range("H74").Value = 2303.5
$ws.Range("I74").Value = 1866.875
$ws.Range("K74").Value = 1866.875
$ws.Range("M74").Value = -992.875

$ws.Range("H77").Value = 2303.5
$ws.Range("I77").Value = 1866.875
$ws.Range("K77").Value = 9334.375
$ws.Range("M77").Value = -4966.375

$ws.Range("H122").Value = 3388.7144
$ws.Range("I122").Value = 2784.5
$ws.Range("J122").Value = 7014
$ws.Range("K122").Value = 8353.5
$ws.Range("L122").Value = 21042
$ws.Range("M122").Value = -5903.5
$ws.Range("N122").Value = -25942

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 225.42857
$ws.Range("I7").Value = 197.25
$ws.Range("K7").Value = 197.25
$ws.Range("M7").Value = -84.25

$ws.Range("H22").Value = 463
$ws.Range("I22").Value = 463
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 463
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -113
$ws.Range("N22").ClearContents()

$ws.Range("H86").Value = 7247.273
$ws.Range("I86").Value = 6185.25
$ws.Range("J86").Value = 7854.143
$ws.Range("K86").Value = 6185.25
$ws.Range("L86").Value = 7854.143
$ws.Range("M86").Value = -5062.25
$ws.Range("N86").Value = -10100.143

$ws.Range("H89").Value = 7247.273
$ws.Range("I89").Value = 6185.25
$ws.Range("J89").Value = 7854.143
$ws.Range("K89").Value = 30926.25
$ws.Range("L89").Value = 39270.715
$ws.Range("M89").Value = -25310.25
$ws.Range("N89").Value = -50502.715

$ws.Range("H107").Value = 1222.6364
$ws.Range("I107").Value = 1233.2858
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 1233.2858
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 686.7141999999999
$ws.Range("N107").Value = -4839

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43997550
$ws.Range("J4").Value = 8841.5
$ws.Range("L4").Value = 26524.5
$ws.Range("N4").Value = -26748.5

$ws.Range("H5").Value = 1658.2858
$ws.Range("I5").Value = 1793.1666
$ws.Range("J5").Value = 849
$ws.Range("K5").Value = 5379.4998
$ws.Range("L5").Value = 2547
$ws.Range("M5").Value = -5267.4998
$ws.Range("N5").Value = -2771

$ws.Range("H14").Value = 546
$ws.Range("I14").Value = 546
$ws.Range("K14").Value = 1638
$ws.Range("M14").Value = -1465

$ws.Range("H22").Value = 5070
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 5070
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

$ws.Range("H29").Value = 7400
$ws.Range("J29").Value = 10850
$ws.Range("L29").Value = 32550
$ws.Range("N29").Value = -33104

$ws.Range("H34").Value = 3856.1667
$ws.Range("I34").Value = 999
$ws.Range("K34").Value = 2997
$ws.Range("M34").Value = -2913

$ws.Range("H40").Value = 553.5
$ws.Range("I40").Value = 43
$ws.Range("J40").Value = 859.8
$ws.Range("K40").Value = 172
$ws.Range("L40").Value = 3439.2
$ws.Range("M40").Value = -103
$ws.Range("N40").Value = -3577.2

$ws.Range("H44").Value = 615.25
$ws.Range("I44").Value = 688.4
$ws.Range("J44").Value = 493.33334
$ws.Range("K44").Value = 2065.2
$ws.Range("L44").Value = 1480.00002
$ws.Range("M44").Value = -1667.2
$ws.Range("N44").Value = -2276.00002

$ws.Range("H86").Value = 702.3333
$ws.Range("I86").Value = 699
$ws.Range("K86").Value = 2097
$ws.Range("M86").Value = -911

$ws.Range("H89").Value = 702.3333
$ws.Range("I89").Value = 699
$ws.Range("K89").Value = 6291
$ws.Range("M89").Value = -363

$ws.Range("H122").Value = 12333
$ws.Range("J122").Value = 12333
$ws.Range("L122").Value = 110997
$ws.Range("N122").Value = -115897

$ws.Range("H129").Value = 764.75
$ws.Range("I129").Value = 764.75
$ws.Range("K129").Value = 2294.25
$ws.Range("M129").Value = 2705.75

$ws.Range("H135").Value = 1658.2858
$ws.Range("I135").Value = 1793.1666
$ws.Range("J135").Value = 849
$ws.Range("K135").Value = 16138.4994
$ws.Range("L135").Value = 7641
$ws.Range("M135").Value = -13603.4994
$ws.Range("N135").Value = -12711

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 4200094
$ws.Range("I14").Value = 5250067.5
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 5250067.5
$ws.Range("L14").Value = 200
$ws.Range("M14").Value = -5249899.5
$ws.Range("N14").Value = -536

$ws.Range("H107").Value = 384.53845
$ws.Range("I107").Value = 336.15384
$ws.Range("K107").Value = 336.15384
$ws.Range("M107").Value = 1583.84616

$ws.Range("H122").Value = 3766.16
$ws.Range("I122").Value = 2611.5
$ws.Range("K122").Value = 7834.5
$ws.Range("M122").Value = -5384.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4309.625
$ws.Range("J46").Value = 3239.182
$ws.Range("L46").Value = 3239.182
$ws.Range("N46").Value = -3615.182

$ws.Range("H95").Value = 79000
$ws.Range("J95").Value = 79000
$ws.Range("L95").Value = 79000
$ws.Range("N95").Value = -84492

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H140").Value = 93484.11
$ws.Range("J140").Value = 88300.875
$ws.Range("L140").Value = 88300.875
$ws.Range("N140").Value = -98660.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 9166.666999999999
$ws.Range("I6").Value = 6250
$ws.Range("J6").Value = 15000
$ws.Range("K6").Value = 6250
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = -6135
$ws.Range("N6").Value = -15230

$ws.Range("H81").Value = 65428.145
$ws.Range("J81").Value = 87799.8
$ws.Range("L81").Value = 175599.6
$ws.Range("N81").Value = -177721.6

$ws.Range("H84").Value = 65428.145
$ws.Range("J84").Value = 87799.8
$ws.Range("L84").Value = 877998
$ws.Range("N84").Value = -888606

$ws.Range("H86").Value = 43999.8
$ws.Range("J86").Value = 43999.8
$ws.Range("L86").Value = 43999.8
$ws.Range("N86").Value = -46245.8

$ws.Range("H89").Value = 43999.8
$ws.Range("J89").Value = 43999.8
$ws.Range("L89").Value = 219999
$ws.Range("N89").Value = -231231

$ws.Range("H132").Value = 3627.2666
$ws.Range("I132").Value = 3775
$ws.Range("K132").Value = 11325
$ws.Range("M132").Value = -8795
